$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "66.748.52"
$ws.Range("E2").Value = "  -4.14%  "
$ws.Range("D3").Value = "3.456.45"
$ws.Range("E3").Value = "  -4.40%  "
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "602.85"
$ws.Range("E5").Value = "  -4.58%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "147.36"
$ws.Range("E6").Value = "  -7.26%  "
$ws.Range("D7").Value = "3.454.76"
$ws.Range("E7").Value = "  -4.39%  "
$ws.Range("E8").Value = "  +0.07%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.483"
$ws.Range("E9").Value = "  -2.41%  "
$ws.Range("E10").Value = "  -5.30%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "7.48"
$ws.Range("E11").Value = "  -0.86%  "
$ws.Range("E12").Value = "  -4.22%  "
$ws.Range("E13").Value = "  -5.75%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "31.65"
$ws.Range("E14").Value = "  -6.27%  "
$ws.Range("D15").Value = "4.043.99"
$ws.Range("E15").Value = "  -4.32%  "
$ws.Range("D16").Value = "3.459.37"
$ws.Range("E16").Value = "  -4.43%  "
$ws.Range("D17").Value = "66.836.78"
$ws.Range("E17").Value = "  -3.83%  "
$ws.Range("E18").Value = "  -1.02%  "
$ws.Range("E19").Value = "  -4.65%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "15.30"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "9.99"
$ws.Range("E21").Value = "  -3.10%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "439.27"
$ws.Range("E22").Value = "  -5.15%  "
$ws.Range("E23").Value = "  -6.10%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "78.60"
$ws.Range("E24").Value = "  -0.55%  "
$ws.Range("E25").Value = "  -0.10%  "
$ws.Range("D26").Value = "3.598.13"
$ws.Range("E26").Value = "  -4.28%  "
$ws.Range("E27").Value = "  -10.71%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.89"
$ws.Range("E28").Value = "  -8.36%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "8.41"
$ws.Range("E29").Value = "  -11.07%  "
$ws.Range("E30").Value = "  -6.70%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.60"
$ws.Range("E31").Value = "  -7.70%  "
$ws.Range("E32").Value = "  -3.06%  "
$ws.Range("E33").Value = "  +0.14%  "
$ws.Range("E34").Value = "  -4.60%  "
$ws.Range("E35").Value = "  -7.87%  "
$ws.Range("D36").Value = "3.450.70"
$ws.Range("E36").Value = "  -4.49%  "
$ws.Range("E37").Value = "  -8.06%  "
$ws.Range("E38").Value = "  +0.00%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "7.90"
$ws.Range("E39").Value = "  -6.87%  "
$ws.Range("E40").Value = "  +0.02%  "
$ws.Range("B41").Value = "Monero"
$ws.Range("C41").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "172.27"
$ws.Range("E41").Value = "  -3.71%  "
$ws.Range("E42").Value = "  -4.47%  "
$ws.Range("B43").Value = "Stacks"
$ws.Range("C43").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.16"
$ws.Range("E43").Value = "  -10.78%  "
$ws.Range("E44").Value = "  -5.58%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.882"
$ws.Range("E45").Value = "  -3.51%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "28.86"
$ws.Range("E46").Value = "  -10.23%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "45.88"
$ws.Range("E47").Value = "  -0.26%  "
$ws.Range("E48").Value = "  -11.04%  "
$ws.Range("E49").Value = "  -4.88%  "
$ws.Range("E50").Value = "  -11.55%  "
$ws.Range("E51").Value = "  -5.85%  "
